$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 ("I0") and J1 ("IF"), styled like the existing header row (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-28 for columns I and J
$data = @{
    2  = @(1, 3)
    3  = @(1, 5)
    4  = @(1, 6)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 7)
    9  = @(1, 5)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 7)
    13 = @(1, 3)
    14 = @(1, 5)
    15 = @(1, 6)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 5)
    19 = @(1, 5)
    20 = @(1, 5)
    21 = @(1, 4)
    22 = @(1, 3)
    23 = @(1, 5)
    24 = @(1, 7)
    25 = @(1, 5)
    26 = @(7, 8)
    27 = @(7, 8)
    28 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
